# The "shard" sheet holds a sorted table (Table1, A1:D40) whose rows are
# kept in ascending order of column D ("UTC"). Row 18 ("hank ulator") had
# the wrong UTC value (16 instead of 17), which put it out of sorted order.
# Fix the value and re-sort the table range so the row lands back in its
# correct sorted position (it ends up at row 24); everything else shifts
# up by one row to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mis-sorted UTC value for "hank ulator".
$ws.Range("D18").Value = 17

# Re-sort the table data (A1:D40, header in row 1) ascending by column D.
$dataRange = $ws.Range("A1:D40")
$sortKey = $ws.Range("D1")
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 1) | Out-Null

# Reflect the resulting selection, as seen in the saved workbook.
$ws.Range("D18").Select() | Out-Null
